$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Standard Excel PasteSpecial constants
$xlPasteValues = -4163
$xlPasteFormats = -4122

# --- Header text updates (rich-text cells: edit only the changed sub-run) ---
# A8 contains "Volume 32   Number  51" -> "...52"
$ws.Cells.Item(8,1).Characters(21,2).Text = "52"

# C9 contains "Report Covering the Week  12/15/2025  Through  12/21/2025"
$ws.Cells.Item(9,3).Characters(27,10).Text = "12/22/2025"
$ws.Cells.Item(9,3).Characters(48,10).Text = "12/28/2025"

# --- Row 15 (Rape) ---
# Several cells collapse from a numeric report into the "no activity" text
# markers ("0" / "***.*") already used elsewhere in the sheet (columns D/E).
# Re-use those donor cells' value + format so the shared-string index and
# cell style exactly match the other "0"/"***.*" cells in the sheet.
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial($xlPasteValues)
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial($xlPasteFormats)

$ws.Range("D15").Copy()
$ws.Range("G15").PasteSpecial($xlPasteValues)
$ws.Range("D15").Copy()
$ws.Range("G15").PasteSpecial($xlPasteFormats)

$ws.Range("E15").Copy()
$ws.Range("H15").PasteSpecial($xlPasteValues)
$ws.Range("E15").Copy()
$ws.Range("H15").PasteSpecial($xlPasteFormats)

$ws.Range("L15").Value = 75

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -83.333333333333
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -57.894736842105
$ws.Range("I16").Value = 189
$ws.Range("J16").Value = 181
$ws.Range("K16").Value = 4.419889502762
$ws.Range("L16").Value = 9.248554913294
$ws.Range("M16").Value = 0.531914893617
$ws.Range("N16").Value = -81.896551724137

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 7.692307692307
$ws.Range("I17").Value = 280
$ws.Range("J17").Value = 229
$ws.Range("K17").Value = 22.270742358078
$ws.Range("L17").Value = 24.444444444444
$ws.Range("M17").Value = 97.183098591549
$ws.Range("N17").Value = -6.976744186046

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -55.555555555555
$ws.Range("G18").Value = 41
$ws.Range("H18").Value = -41.463414634146
$ws.Range("I18").Value = 376
$ws.Range("J18").Value = 253
$ws.Range("K18").Value = 48.616600790513
$ws.Range("L18").Value = 72.477064220183
$ws.Range("M18").Value = 48.616600790513
$ws.Range("N18").Value = -77.101096224116

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -26.315789473684
$ws.Range("F19").Value = 85
$ws.Range("G19").Value = 94
$ws.Range("H19").Value = -9.574468085106
$ws.Range("I19").Value = 1239
$ws.Range("J19").Value = 1086
$ws.Range("K19").Value = 14.088397790055
$ws.Range("L19").Value = 16.447368421052
$ws.Range("M19").Value = -15.426621160409
$ws.Range("N19").Value = -54.863387978142

# --- Row 20 (G.L.A.) ---
$ws.Range("D20").Copy()
$ws.Range("C20").PasteSpecial($xlPasteValues)
$ws.Range("D20").Copy()
$ws.Range("C20").PasteSpecial($xlPasteFormats)

$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 400
$ws.Range("L20").Value = -19.696969696969
$ws.Range("N20").Value = -95.839874411303

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -34.285714285714
$ws.Range("F21").Value = 139
$ws.Range("G21").Value = 168
$ws.Range("H21").Value = -17.261904761904
$ws.Range("I21").Value = 2159
$ws.Range("J21").Value = 1814
$ws.Range("K21").Value = 19.018743109151
$ws.Range("L21").Value = 22.670454545454
$ws.Range("M21").Value = 1.935788479697
$ws.Range("N21").Value = -69.301862647518

# --- Row 22 (Transit) ---
$ws.Range("D22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -71.428571428571
$ws.Range("J22").Value = 65
$ws.Range("K22").Value = 20
$ws.Range("L22").Value = -13.333333333333
$ws.Range("M22").Value = 1.298701298701

# --- Row 23 (Housing) ---
$ws.Range("F23").Value = 2
$ws.Range("H23").Value = 100

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 50
$ws.Range("D24").Value = 51
$ws.Range("E24").Value = -1.960784313725
$ws.Range("F24").Value = 257
$ws.Range("G24").Value = 251
$ws.Range("H24").Value = 2.390438247011
$ws.Range("I24").Value = 2888
$ws.Range("J24").Value = 3065
$ws.Range("K24").Value = -5.774877650897
$ws.Range("L24").Value = 27.056753189617
$ws.Range("M24").Value = 57.298474945533

# --- Row 25 (Retail Theft) ---
$ws.Range("C25").Value = 40
$ws.Range("D25").Value = 47
$ws.Range("E25").Value = -14.893617021276
$ws.Range("F25").Value = 202
$ws.Range("G25").Value = 227
$ws.Range("H25").Value = -11.013215859030
$ws.Range("I25").Value = 2386
$ws.Range("J25").Value = 2626
$ws.Range("K25").Value = -9.139375476009
$ws.Range("L25").Value = 40.933254577672

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 46
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = 43.75
$ws.Range("I26").Value = 592
$ws.Range("J26").Value = 605
$ws.Range("K26").Value = -2.148760330578
$ws.Range("L26").Value = 15.175097276264
$ws.Range("M26").Value = 43.689320388349

# --- Row 27 (UCR Rape*) ---
$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial($xlPasteValues)
$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial($xlPasteFormats)

$ws.Range("D27").Copy()
$ws.Range("G27").PasteSpecial($xlPasteValues)
$ws.Range("D27").Copy()
$ws.Range("G27").PasteSpecial($xlPasteFormats)

$ws.Range("E27").Copy()
$ws.Range("H27").PasteSpecial($xlPasteValues)
$ws.Range("E27").Copy()
$ws.Range("H27").PasteSpecial($xlPasteFormats)

$ws.Range("L27").Value = 4.166666666666

# --- Row 28 (Other Sex Crimes) ---
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = -50
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 122
$ws.Range("J28").Value = 121
$ws.Range("K28").Value = 0.826446280991
$ws.Range("L28").Value = 9.909909909909

# --- Row 31 (Hate Crimes) ---
$ws.Range("I31").Value = 9
$ws.Range("K31").Value = -35.714285714285
$ws.Range("L31").Value = -35.714285714285

$excel.CutCopyMode = $false
